$d = $word.ActiveDocument

# Helper: locate the paragraph index (1-based) that fully contains a given Range.
function Get-ParaIndexForRange($rng) {
    $idx = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            $idx = $i
        }
    }
    return $idx
}

# ---------------------------------------------------------------------------
# 1) "List ToDo in Code" heading: mark the implicit "ToDo in" phrase as a
#    grammar-flagged span and split the " in " run into " in" + " " so a
#    <w:proofErr w:type="gramEnd"/> boundary can sit between them.
# ---------------------------------------------------------------------------
$headingHit = $d.Content
$headingHit.Find.Execute(" in ")
$inStart = $headingHit.Start
$inRange = $d.Range($inStart, $inStart + 4)
$inRange.Text = " in"

$spacePos = $inStart + 3
$spaceIns = $d.Range($spacePos, $spacePos)
$spaceIns.InsertAfter(" ")

# Force the autocoalesced text apart into two distinct runs (" in" / " ")
# by toggling a direct-formatting property off again.
$splitRange = $d.Range($spacePos, $spacePos + 1)
$splitRange.Bold = 1
$splitRange.Bold = 0

# ---------------------------------------------------------------------------
# 2) "Login usuarios" bullet becomes bold (active form).
# ---------------------------------------------------------------------------
$loginHit = $d.Content
$loginHit.Find.Execute("Login usuarios")
$loginIdx = Get-ParaIndexForRange $loginHit
$loginPara = $d.Paragraphs.Item($loginIdx)
$loginPara.Range.Bold = 1
$loginPara.Range.BoldBi = 1

# ---------------------------------------------------------------------------
# 3) Insert a new "Sign in usuarios" bullet right after "Login usuarios",
#    cloned from it (so it inherits the same list/paragraph formatting),
#    then retext it and bold it.
# ---------------------------------------------------------------------------
$nextPara = $d.Paragraphs.Item($loginIdx + 1)
$cloneSrc = $d.Range($loginPara.Range.Start, $nextPara.Range.Start)
$cloneDst = $d.Range($nextPara.Range.Start, $nextPara.Range.Start)
$cloneDst.FormattedText = $cloneSrc.FormattedText

$signPara = $d.Paragraphs.Item($loginIdx + 1)
$signStart = $signPara.Range.Start

# "Login" -> "Sign"
$wordRange = $d.Range($signStart, $signStart + 5)
$wordRange.Text = "Sign"

# " usuarios" -> " in usuarios"
$afterSign = $signStart + 4
$insPoint = $d.Range($afterSign, $afterSign)
$insPoint.InsertBefore(" in")

# Split "Sign" and " in usuarios" into separate runs.
$splitRange2 = $d.Range($signStart, $afterSign)
$splitRange2.Bold = 1
$splitRange2.Bold = 0

$signPara.Range.Bold = 1
$signPara.Range.BoldBi = 1

# ---------------------------------------------------------------------------
# 4) "Cookies de Login" bullet becomes bold (active form) as well.
# ---------------------------------------------------------------------------
$cookiesHit = $d.Content
$cookiesHit.Find.Execute("Cookies de Login")
$cookiesIdx = Get-ParaIndexForRange $cookiesHit
$cookiesPara = $d.Paragraphs.Item($cookiesIdx)
$cookiesPara.Range.Bold = 1
$cookiesPara.Range.BoldBi = 1

Write-Output "done"
